$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new text looks like a plain number (e.g. "1.011")
# must be forced to Text format first, otherwise Excel auto-converts them
# to a numeric value and drops formatting such as trailing zeros.

$ws.Range("D2").Value = "29.555.01"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.907.22"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.24"
$ws.Range("E5").Value = "  +3.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4004"
$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08041"
$ws.Range("E9").Value = "  -2.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9887"
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("D12").Value = "1.915.02"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.923"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.118"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.20"
$ws.Range("E15").Value = "  -2.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06834"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001020"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.38"
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "29.570.60"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.517"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("E23").Value = "  -1.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.153"
$ws.Range("E24").Value = "  -1.16%  "

$ws.Range("D25").Value = "2.157.78"
$ws.Range("E25").Value = "  +1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.90"
$ws.Range("E26").Value = "  +0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.524"
$ws.Range("E27").Value = "  -2.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.67"
$ws.Range("E28").Value = "  -1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.055"
$ws.Range("E29").Value = "  -2.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.23"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9935"
$ws.Range("E31").Value = "  -2.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09514"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.494"
$ws.Range("E33").Value = "  -3.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.551"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.389"
$ws.Range("E35").Value = "  +1.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06457"
$ws.Range("E36").Value = "  +5.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02245"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.197"
$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5823"
$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.54"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.767"
$ws.Range("E41").Value = "  -4.01%  "

$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.442"
$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.273"
$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07430"
$ws.Range("E45").Value = "  -2.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.16"
$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5472"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.938"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.13"
$ws.Range("E49").Value = "  -1.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.380"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.30"
$ws.Range("E51").Value = "  -1.61%  "
